$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "MCT-2A-Sistemas digitais"
$ws.Range("B4").Value = "MCT-2A-Sistemas digitais"
$ws.Range("C4").Value = "-"
$ws.Range("D4").Value = "ELT-2A-Sistemas digitais"
$ws.Range("B6").Value = "-"
$ws.Range("D6").Value = "ELT-2A-Sistemas digitais"
$ws.Range("B7").Value = "-"
$ws.Range("F7").Value = "-"
